$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2345
$ws.Range("I38").Value = 175.81818
$ws.Range("J38").Value = 4996.222
$ws.Range("K38").Value = 527.4545400000001
$ws.Range("L38").Value = 14988.666
$ws.Range("M38").Value = -155.4545400000001
$ws.Range("N38").Value = -15732.666
$ws.Range("H40").Value = 4705.75
$ws.Range("I40").Value = 3229
$ws.Range("J40").Value = 5198
$ws.Range("K40").Value = 3229
$ws.Range("L40").Value = 5198
$ws.Range("M40").Value = -3054
$ws.Range("N40").Value = -5548
$ws.Range("H43").Value = 4946.7207
$ws.Range("I43").Value = 5477.2
$ws.Range("J43").Value = 4209.9443
$ws.Range("K43").Value = 5477.2
$ws.Range("L43").Value = 4209.9443
$ws.Range("M43").Value = -5408.2
$ws.Range("N43").Value = -4347.9443
$ws.Range("H51").Value = 7786.1816
$ws.Range("I51").Value = 5749.75
$ws.Range("J51").Value = 8949.857
$ws.Range("K51").Value = 5749.75
$ws.Range("L51").Value = 8949.857
$ws.Range("M51").Value = -5265.75
$ws.Range("N51").Value = -9917.857
$ws.Range("H101").Value = 644.5
$ws.Range("J101").Value = 534
$ws.Range("L101").Value = 1602
$ws.Range("N101").Value = -4846
$ws.Range("H116").Value = 12375
$ws.Range("I116").Value = 11747
$ws.Range("J116").Value = 13003
$ws.Range("K116").Value = 11747
$ws.Range("L116").Value = 13003
$ws.Range("M116").Value = -8305
$ws.Range("N116").Value = -19887
$ws.Range("H138").Value = 5452.486
$ws.Range("I138").Value = 4020.625
$ws.Range("J138").Value = 6658.263
$ws.Range("K138").Value = 12061.875
$ws.Range("L138").Value = 19974.789
$ws.Range("M138").Value = -6921.875
$ws.Range("N138").Value = -30254.789

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 289.6
$ws.Range("I5").Value = 274.5
$ws.Range("J5").Value = 350
$ws.Range("K5").Value = 274.5
$ws.Range("L5").Value = 350
$ws.Range("M5").Value = -162.5
$ws.Range("N5").Value = -574
$ws.Range("H88").Value = 1117.5
$ws.Range("J88").Value = 1144.9166
$ws.Range("L88").Value = 1144.9166
$ws.Range("N88").Value = -1956.9166
$ws.Range("H91").Value = 1117.5
$ws.Range("J91").Value = 1144.9166
$ws.Range("L91").Value = 1144.9166
$ws.Range("N91").Value = -3952.9166
$ws.Range("H106").Value = 149645.42
$ws.Range("J106").Value = 149645.42
$ws.Range("L106").Value = 149645.42
$ws.Range("N106").Value = -152169.42
$ws.Range("H122").Value = 4613.36
$ws.Range("I122").Value = 4671.4546
$ws.Range("J122").Value = 4187.3335
$ws.Range("K122").Value = 14014.3638
$ws.Range("L122").Value = 12562.0005
$ws.Range("M122").Value = -11564.3638
$ws.Range("N122").Value = -17462.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 289.6
$ws.Range("I4").Value = 274.5
$ws.Range("J4").Value = 350
$ws.Range("K4").Value = 274.5
$ws.Range("L4").Value = 350
$ws.Range("M4").Value = -159.5
$ws.Range("N4").Value = -580
$ws.Range("H22").Value = 468.5
$ws.Range("I22").Value = 416.85715
$ws.Range("K22").Value = 416.85715
$ws.Range("M22").Value = -243.85715
$ws.Range("H80").Value = 29271.143
$ws.Range("I80").Value = 33963.668
$ws.Range("J80").Value = 25751.75
$ws.Range("K80").Value = 33963.668
$ws.Range("L80").Value = 25751.75
$ws.Range("M80").Value = -32965.668
$ws.Range("N80").Value = -27747.75
$ws.Range("H83").Value = 29271.143
$ws.Range("I83").Value = 33963.668
$ws.Range("J83").Value = 25751.75
$ws.Range("K83").Value = 169818.34
$ws.Range("L83").Value = 128758.75
$ws.Range("M83").Value = -164826.34
$ws.Range("N83").Value = -138742.75
$ws.Range("H94").Value = 832.5909
$ws.Range("I94").Value = 548.2941
$ws.Range("J94").Value = 1799.2
$ws.Range("K94").Value = 548.2941
$ws.Range("L94").Value = 1799.2
$ws.Range("M94").Value = -97.29409999999996
$ws.Range("N94").Value = -2701.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11674.929
$ws.Range("I31").Value = 5268.091
$ws.Range("K31").Value = 5268.091
$ws.Range("M31").Value = -4973.091
$ws.Range("H34").Value = 11674.929
$ws.Range("I34").Value = 5268.091
$ws.Range("K34").Value = 5268.091
$ws.Range("M34").Value = -5066.091
$ws.Range("H58").Value = 4949.788
$ws.Range("I58").Value = 1765.84
$ws.Range("J58").Value = 14899.625
$ws.Range("K58").Value = 1765.84
$ws.Range("L58").Value = 14899.625
$ws.Range("M58").Value = -1562.84
$ws.Range("N58").Value = -15305.625
$ws.Range("H69").Value = 130494.6
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 130494.6
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 130494.6
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -131992.6
$ws.Range("H72").Value = 130494.6
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 130494.6
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 391483.8
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -398971.8
$ws.Range("H122").Value = 1933.2667
$ws.Range("J122").Value = 3749.5
$ws.Range("L122").Value = 11248.5
$ws.Range("N122").Value = -16148.5
$ws.Range("H135").Value = 91992.5
$ws.Range("J135").Value = 91992.5
$ws.Range("L135").Value = 91992.5
$ws.Range("N135").Value = -102132.5
$ws.Range("H136").Value = 4949.788
$ws.Range("I136").Value = 1765.84
$ws.Range("J136").Value = 14899.625
$ws.Range("K136").Value = 5297.52
$ws.Range("L136").Value = 44698.875
$ws.Range("M136").Value = -2747.52
$ws.Range("N136").Value = -49798.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H129").Value = 16868328
$ws.Range("I129").Value = 168861
$ws.Range("J129").Value = 41917530
$ws.Range("K129").Value = 506583
$ws.Range("L129").Value = 125752590
$ws.Range("M129").Value = -501583
$ws.Range("N129").Value = -125762590

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4173.224
$ws.Range("I70").Value = 2797.7144
$ws.Range("J70").Value = 4801.174
$ws.Range("K70").Value = 2797.7144
$ws.Range("L70").Value = 4801.174
$ws.Range("M70").Value = -2527.7144
$ws.Range("N70").Value = -5341.174
$ws.Range("H73").Value = 4173.224
$ws.Range("I73").Value = 2797.7144
$ws.Range("J73").Value = 4801.174
$ws.Range("K73").Value = 2797.7144
$ws.Range("L73").Value = 4801.174
$ws.Range("M73").Value = -1861.7144
$ws.Range("N73").Value = -6673.174
$ws.Range("H102").Value = 3626.7856
$ws.Range("I102").Value = 1749.9
$ws.Range("J102").Value = 8319
$ws.Range("K102").Value = 1749.9
$ws.Range("L102").Value = 8319
$ws.Range("M102").Value = -127.9000000000001
$ws.Range("N102").Value = -11563
$ws.Range("H122").Value = 5632.759
$ws.Range("I122").Value = 4333
$ws.Range("K122").Value = 12999
$ws.Range("M122").Value = -10549

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 814.4091
$ws.Range("I22").Value = 344
$ws.Range("K22").Value = 344
$ws.Range("M22").Value = -49
$ws.Range("H27").Value = 814.4091
$ws.Range("I27").Value = 344
$ws.Range("K27").Value = 344
$ws.Range("M27").Value = -237
$ws.Range("H46").Value = 6071.4287
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 6500
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 6500
$ws.Range("M46").Value = -4812
$ws.Range("N46").Value = -6876
$ws.Range("H55").Value = 1096.5264
$ws.Range("J55").Value = 1248.2307
$ws.Range("L55").Value = 1248.2307
$ws.Range("N55").Value = -1594.2307
$ws.Range("H62").Value = 49999
$ws.Range("J62").Value = 49999
$ws.Range("L62").Value = 49999
$ws.Range("N62").Value = -51247
$ws.Range("H65").Value = 49999
$ws.Range("J65").Value = 49999
$ws.Range("L65").Value = 149997
$ws.Range("N65").Value = -156237

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 116526
$ws.Range("J106").Value = 116526
$ws.Range("L106").Value = 116526
$ws.Range("N106").Value = -119050
$ws.Range("H132").Value = 3587.1875
$ws.Range("I132").Value = 1400.1818
$ws.Range("J132").Value = 8398.6
$ws.Range("K132").Value = 4200.5454
$ws.Range("L132").Value = 25195.8
$ws.Range("M132").Value = -1670.5454
$ws.Range("N132").Value = -30255.8
